$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Cells.Item(2, 4).Value = "34.070.20"
$ws.Cells.Item(2, 5).Value = "  -1.76%  "

# Row 3 (Ethereum)
$ws.Cells.Item(3, 4).Value = "1.787.07"
$ws.Cells.Item(3, 5).Value = "  -1.20%  "

# Row 4 (TetherUSD)
$ws.Cells.Item(4, 5).Value = "  +0.13%  "

# Row 5 (BNB)
$ws.Cells.Item(5, 4).Value = "222.86"
$ws.Cells.Item(5, 5).Value = "  -1.19%  "

# Row 7 (USDC)
$ws.Cells.Item(7, 5).Value = "  +0.11%  "

# Row 8 (Solana)
$ws.Cells.Item(8, 4).Value = "32.24"
$ws.Cells.Item(8, 5).Value = "  -1.38%  "

# Row 9 (Cardano)
$ws.Cells.Item(9, 4).Value = "0.284"
$ws.Cells.Item(9, 5).Value = "  -2.11%  "

# Row 10 (Dogecoin)
$ws.Cells.Item(10, 5).Value = "  -0.39%  "

# Row 11 (TRON)
$ws.Cells.Item(11, 4).Value = "0.0930"
$ws.Cells.Item(11, 5).Value = "  +0.14%  "

# Row 12 (WrappedliquidstakedEther2.0)
$ws.Cells.Item(12, 4).Value = "2.044.65"
$ws.Cells.Item(12, 5).Value = "  -1.15%  "

# Row 13 (Chainlink)
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.794.58"
$ws.Cells.Item(13, 5).Value = "  -0.72%  "

# Row 14 (WrappedEther)
$ws.Cells.Item(14, 2).Value = "Chainlink"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(14, 4).Value = "10.94"
$ws.Cells.Item(14, 5).Value = "  -1.72%  "

# Row 15 (Polygon)
$ws.Cells.Item(15, 4).Value = "0.624"
$ws.Cells.Item(15, 5).Value = "  -3.23%  "

# Row 16 (WrappedBTC)
$ws.Cells.Item(16, 4).Value = "34.056.17"
$ws.Cells.Item(16, 5).Value = "  -1.82%  "

# Row 17 (Polkadot)
$ws.Cells.Item(17, 4).Value = "4.16"
$ws.Cells.Item(17, 5).Value = "  -4.04%  "

# Row 18 (Litecoin)
$ws.Cells.Item(18, 4).Value = "67.91"
$ws.Cells.Item(18, 5).Value = "  -2.58%  "

# Row 19 (BitcoinCash)
$ws.Cells.Item(19, 4).Value = "243.86"
$ws.Cells.Item(19, 5).Value = "  -4.24%  "

# Row 20 (ShibaInu)
$ws.Cells.Item(20, 4).Value = "0.0₃0782"
$ws.Cells.Item(20, 5).Value = "  -2.91%  "

# Row 21 (Dai)
$ws.Cells.Item(21, 5).Value = "  +0.16%  "

# Row 22 (Avalanche)
$ws.Cells.Item(22, 5).Value = "  -1.96%  "

# Row 23 (Uniswap)
$ws.Cells.Item(23, 5).Value = "  -4.46%  "

# Row 24 (Toncoin)
$ws.Cells.Item(24, 4).Value = "2.12"
$ws.Cells.Item(24, 5).Value = "  -2.43%  "

# Row 25 (Monero)
$ws.Cells.Item(25, 4).Value = "158.64"
$ws.Cells.Item(25, 5).Value = "  -1.91%  "

# Row 26 (EthereumClassic)
$ws.Cells.Item(26, 4).Value = "16.28"
$ws.Cells.Item(26, 5).Value = "  -1.39%  "

# Row 27 (Cosmos)
$ws.Cells.Item(27, 5).Value = "  -2.16%  "

# Row 28 (Stellar)
$ws.Cells.Item(28, 5).Value = "  -2.24%  "

# Row 29 (BinanceUSD)
$ws.Cells.Item(29, 5).Value = "  +0.28%  "

# Row 30 (Hedera)
$ws.Cells.Item(30, 4).Value = "0.0518"
$ws.Cells.Item(30, 5).Value = "  -2.96%  "

# Row 31 (PancakeSwap)
$ws.Cells.Item(31, 4).Value = "1.20"
$ws.Cells.Item(31, 5).Value = "  -0.18%  "

# Row 32 (Filecoin)
$ws.Cells.Item(32, 5).Value = "  -3.98%  "

# Row 33 (InternetComputer(DFINITY))
$ws.Cells.Item(33, 5).Value = "  -4.45%  "

# Row 34 (LidoDAOToken)
$ws.Cells.Item(34, 5).Value = "  -5.23%  "

# Row 35 (Maker)
$ws.Cells.Item(35, 4).Value = "1.384.47"
$ws.Cells.Item(35, 5).Value = "  -3.89%  "

# Row 36 (ImmutableX)
$ws.Cells.Item(36, 4).Value = "0.646"
$ws.Cells.Item(36, 5).Value = "  +0.16%  "

# Row 37 (TrustWalletToken)
$ws.Cells.Item(37, 5).Value = "  -1.95%  "

# Row 38 (VeChain)
$ws.Cells.Item(38, 4).Value = "0.0184"
$ws.Cells.Item(38, 5).Value = "  -4.31%  "

# Row 39 (Aave)
$ws.Cells.Item(39, 4).Value = "79.53"
$ws.Cells.Item(39, 5).Value = "  -6.62%  "

# Row 40 (HuobiToken)
$ws.Cells.Item(40, 5).Value = "  +0.34%  "

# Row 41 (ARBITRUM)
$ws.Cells.Item(41, 4).Value = "0.914"
$ws.Cells.Item(41, 5).Value = "  -4.74%  "

# Row 42 (MXToken)
$ws.Cells.Item(42, 5).Value = "  -3.86%  "

# Row 43 (RenderToken)
$ws.Cells.Item(43, 4).Value = "2.16"
$ws.Cells.Item(43, 5).Value = "  +0.12%  "

# Row 44 (Kaspa)
$ws.Cells.Item(44, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(44, 4).Value = "0.0₆0136"
$ws.Cells.Item(44, 5).Value = "  +7.06%  "

# Row 45 (WEMIXToken)
$ws.Cells.Item(45, 2).Value = "Kaspa"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(45, 4).Value = "0.0498"
$ws.Cells.Item(45, 5).Value = "  +0.88%  "

# Row 46 (FraxShare)
$ws.Cells.Item(46, 2).Value = "WEMIXToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(46, 4).Value = "1.05"
$ws.Cells.Item(46, 5).Value = "  -0.73%  "

# Row 47 (Quant)
$ws.Cells.Item(47, 4).Value = "107.24"
$ws.Cells.Item(47, 5).Value = "  +0.85%  "

# Row 48 (BabyDogeCoin)
$ws.Cells.Item(48, 2).Value = "FraxShare"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(48, 4).Value = "5.86"
$ws.Cells.Item(48, 5).Value = "  -3.25%  "

# Row 49 (RocketPoolETH)
$ws.Cells.Item(49, 4).Value = "1.943.65"
$ws.Cells.Item(49, 5).Value = "  -0.84%  "

# Row 50 (PaxDollar)
$ws.Cells.Item(50, 5).Value = "  +0.00%  "

# Row 51 (InjectiveProtocol)
$ws.Cells.Item(51, 4).Value = "11.98"
$ws.Cells.Item(51, 5).Value = "  -2.07%  "
